$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2315.389
$ws.Range("I15").Value = 2315.389
$ws.Range("K15").Value = 6946.167
$ws.Range("M15").Value = -6777.167

$ws.Range("H33").Value = 220.57143
$ws.Range("I33").Value = 230.46153
$ws.Range("J33").Value = 92
$ws.Range("K33").Value = 230.46153
$ws.Range("L33").Value = 92
$ws.Range("M33").Value = -1.46153000000001
$ws.Range("N33").Value = -550

$ws.Range("H39").Value = 197.42857
$ws.Range("I39").Value = 147
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 441
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -145
$ws.Range("N39").Value = -2092

$ws.Range("H62").Value = 3624.5
$ws.Range("I62").Value = 4499.5
$ws.Range("K62").Value = 4499.5
$ws.Range("M62").Value = -3875.5

$ws.Range("H65").Value = 3624.5
$ws.Range("I65").Value = 4499.5
$ws.Range("K65").Value = 22497.5
$ws.Range("M65").Value = -19377.5

$ws.Range("H98").Value = 603.44446
$ws.Range("I98").Value = 603.44446
$ws.Range("K98").Value = 603.44446
$ws.Range("M98").Value = 894.55554

$ws.Range("H122").Value = 603.44446
$ws.Range("I122").Value = 603.44446
$ws.Range("K122").Value = 1810.33338
$ws.Range("M122").Value = 639.66662

$ws.Range("H137").Value = 3588
$ws.Range("I137").Value = 751.1667
$ws.Range("K137").Value = 2253.5001
$ws.Range("M137").Value = 296.4998999999998

$ws.Range("H138").Value = 4276.8057
$ws.Range("J138").Value = 4462.387
$ws.Range("L138").Value = 13387.161
$ws.Range("N138").Value = -23667.161

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 853.8
$ws.Range("I2").Value = 846.8570999999999
$ws.Range("K2").Value = 846.8570999999999
$ws.Range("M2").Value = -733.8570999999999

$ws.Range("H32").Value = 7574.0713
$ws.Range("I32").Value = 6810.5386
$ws.Range("J32").Value = 17500
$ws.Range("K32").Value = 6810.5386
$ws.Range("L32").Value = 17500
$ws.Range("M32").Value = -6523.5386
$ws.Range("N32").Value = -18074

$ws.Range("H37").Value = 35492.25
$ws.Range("J37").Value = 35492.25
$ws.Range("L37").Value = 35492.25
$ws.Range("N37").Value = -36038.25

$ws.Range("H45").Value = 2334.8125
$ws.Range("I45").Value = 2157.125
$ws.Range("J45").Value = 2512.5
$ws.Range("K45").Value = 2157.125
$ws.Range("L45").Value = 2512.5
$ws.Range("M45").Value = -1780.125
$ws.Range("N45").Value = -3266.5

$ws.Range("H61").Value = 2796
$ws.Range("I61").Value = 2645.75
$ws.Range("K61").Value = 2645.75
$ws.Range("M61").Value = -2433.75

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null

$ws.Range("H97").Value = 551
$ws.Range("I97").Value = 551
$ws.Range("K97").Value = 551
$ws.Range("M97").Value = -55

$ws.Range("H110").Value = 2321
$ws.Range("I110").Value = 2194.7693
$ws.Range("K110").Value = 2194.7693
$ws.Range("M110").Value = -149.7692999999999

$ws.Range("H116").Value = 853.8
$ws.Range("I116").Value = 846.8570999999999
$ws.Range("K116").Value = 846.8570999999999
$ws.Range("M116").Value = 1447.1429

$ws.Range("H122").Value = 6525.091
$ws.Range("I122").Value = 6525.091
$ws.Range("K122").Value = 19575.273
$ws.Range("M122").Value = -17125.273

$ws.Range("H136").Value = 2796
$ws.Range("I136").Value = 2645.75
$ws.Range("K136").Value = 7937.25
$ws.Range("M136").Value = -5387.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 853.8
$ws.Range("I3").Value = 846.8570999999999
$ws.Range("K3").Value = 846.8570999999999
$ws.Range("M3").Value = -732.8570999999999

$ws.Range("H107").Value = 1255.125
$ws.Range("I107").Value = 1548.2
$ws.Range("K107").Value = 1548.2
$ws.Range("M107").Value = 371.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2999
$ws.Range("J16").Value = 2999
$ws.Range("L16").Value = 2999
$ws.Range("N16").Value = -3573

$ws.Range("H113").Value = 2999
$ws.Range("J113").Value = 2999
$ws.Range("L113").Value = 2999
$ws.Range("N113").Value = -7339

$ws.Range("H134").Value = 5999.3335
$ws.Range("I134").Value = 5999.3335
$ws.Range("K134").Value = 17998.0005
$ws.Range("M134").Value = -15463.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19875
$ws.Range("J15").Value = 19875
$ws.Range("L15").Value = 19875
$ws.Range("N15").Value = -20451

$ws.Range("H46").Value = 39946
$ws.Range("J46").Value = 39946
$ws.Range("L46").Value = 39946
$ws.Range("N46").Value = -40258

$ws.Range("H57").Value = 38000
$ws.Range("J57").Value = 38000
$ws.Range("L57").Value = 38000
$ws.Range("N57").Value = -39640

$ws.Range("H80").Value = 3643.7
$ws.Range("I80").Value = 2498.6365
$ws.Range("J80").Value = 5043.222
$ws.Range("K80").Value = 2498.6365
$ws.Range("L80").Value = 5043.222
$ws.Range("M80").Value = -1500.6365
$ws.Range("N80").Value = -7039.222

$ws.Range("H81").Value = 19875
$ws.Range("J81").Value = 19875
$ws.Range("L81").Value = 19875
$ws.Range("N81").Value = -21871

$ws.Range("H83").Value = 3643.7
$ws.Range("I83").Value = 2498.6365
$ws.Range("J83").Value = 5043.222
$ws.Range("K83").Value = 12493.1825
$ws.Range("L83").Value = 25216.11
$ws.Range("M83").Value = -7501.182500000001
$ws.Range("N83").Value = -35200.11

$ws.Range("H84").Value = 19875
$ws.Range("J84").Value = 19875
$ws.Range("L84").Value = 59625
$ws.Range("N84").Value = -69609

$ws.Range("H97").Value = 814.44446
$ws.Range("I97").Value = 666.25
$ws.Range("K97").Value = 666.25
$ws.Range("M97").Value = -170.25

$ws.Range("H122").Value = 972.36365
$ws.Range("I122").Value = 919.7
$ws.Range("K122").Value = 2759.1
$ws.Range("M122").Value = -309.1000000000004

$ws.Range("H126").Value = 2444.25
$ws.Range("J126").Value = 2666.6667
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("N126").Value = -12940.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 439.2857
$ws.Range("I55").Value = 471.25
$ws.Range("J55").Value = 396.66666
$ws.Range("K55").Value = 471.25
$ws.Range("L55").Value = 396.66666
$ws.Range("M55").Value = -298.25
$ws.Range("N55").Value = -742.66666

$ws.Range("H61").Value = 4470
$ws.Range("I61").Value = 4469.857
$ws.Range("K61").Value = 4469.857
$ws.Range("M61").Value = -4267.857

$ws.Range("H113").Value = 4470
$ws.Range("I113").Value = 4469.857
$ws.Range("K113").Value = 4469.857
$ws.Range("M113").Value = -2299.857

$ws.Range("H122").Value = 2579.8
$ws.Range("I122").Value = 2724.75
$ws.Range("K122").Value = 8174.25
$ws.Range("M122").Value = -5724.25

$ws.Range("H132").Value = 4468.8823
$ws.Range("I132").Value = 4100.143
$ws.Range("J132").Value = 4727
$ws.Range("K132").Value = 12300.429
$ws.Range("L132").Value = 14181
$ws.Range("M132").Value = -9770.429
$ws.Range("N132").Value = -19241

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5001.5
$ws.Range("J62").Value = 5003
$ws.Range("L62").Value = 5003
$ws.Range("N62").Value = -6251

$ws.Range("H65").Value = 5001.5
$ws.Range("J65").Value = 5003
$ws.Range("L65").Value = 25015
$ws.Range("N65").Value = -31255

$ws.Range("H113").Value = 745.1429000000001
$ws.Range("I113").Value = 758
$ws.Range("K113").Value = 2274
$ws.Range("M113").Value = -104

$ws.Range("H122").Value = 3959.2
$ws.Range("I122").Value = 3949.25
$ws.Range("K122").Value = 11847.75
$ws.Range("M122").Value = -9397.75

$ws.Range("H132").Value = 3091.3103
$ws.Range("I132").Value = 2806.913
$ws.Range("K132").Value = 8420.739
$ws.Range("M132").Value = -5890.739
